$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_34_7_24"
$ws.Range("B2").Value = [double]"0.9999880372142435"
$ws.Range("C2").Value = [double]"0.9990500038347232"
$ws.Range("D2").Value = [double]"0.9998269074737384"
$ws.Range("E2").Value = [double]"0.9999097463896631"
$ws.Range("F2").Value = [double]"0.9999558291228156"
$ws.Range("G2").Value = [double]"1.116673539815936e-05"
$ws.Range("H2").Value = [double]"0.0008867797202737151"
$ws.Range("I2").Value = [double]"4.320004712944014e-05"
$ws.Range("J2").Value = [double]"2.594710888469202e-05"
$ws.Range("K2").Value = [double]"3.457357800706607e-05"
$ws.Range("L2").Value = [double]"0.0002401447512923727"
$ws.Range("M2").Value = [double]"0.003341666559990592"
$ws.Range("N2").Value = [double]"1.000009900236488"
$ws.Range("O2").Value = [double]"0.003483928223652795"
$ws.Range("P2").Value = [double]"128.8051425055225"
$ws.Range("Q2").Value = [double]"193.4055612235371"

$ws.Range("A3").Value = "model_34_7_23"
$ws.Range("B3").Value = [double]"0.9999882077057538"
$ws.Range("C3").Value = [double]"0.9990494652399772"
$ws.Range("D3").Value = [double]"0.9998297822285379"
$ws.Range("E3").Value = [double]"0.9999112381476529"
$ws.Range("F3").Value = [double]"0.9999565614006991"
$ws.Range("G3").Value = [double]"1.100758905694084e-05"
$ws.Range("H3").Value = [double]"0.0008872824748275157"
$ws.Range("I3").Value = [double]"4.248257223027802e-05"
$ws.Range("J3").Value = [double]"2.551824175299648e-05"
$ws.Range("K3").Value = [double]"3.400040699163725e-05"
$ws.Range("L3").Value = [double]"0.0002407665975201593"
$ws.Range("M3").Value = [double]"0.003317768686473009"
$ws.Range("N3").Value = [double]"1.000009759140066"
$ws.Range("O3").Value = [double]"0.00345901296818415"
$ws.Range("P3").Value = [double]"128.83385121756"
$ws.Range("Q3").Value = [double]"193.4342699355746"

$ws.Range("A4").Value = "model_34_7_22"
$ws.Range("B4").Value = [double]"0.9999883923315521"
$ws.Range("C4").Value = [double]"0.9990488603245904"
$ws.Range("D4").Value = [double]"0.9998329505619905"
$ws.Range("E4").Value = [double]"0.9999128855281826"
$ws.Range("F4").Value = [double]"0.9999573690633251"
$ws.Range("G4").Value = [double]"1.083524897837294e-05"
$ws.Range("H4").Value = [double]"0.0008878471367883882"
$ws.Range("I4").Value = [double]"4.169182662484029e-05"
$ws.Range("J4").Value = [double]"2.504463452755205e-05"
$ws.Range("K4").Value = [double]"3.336823057619617e-05"
$ws.Range("L4").Value = [double]"0.0002414249129454799"
$ws.Range("M4").Value = [double]"0.003291693937530181"
$ws.Range("N4").Value = [double]"1.000009606346302"
$ws.Range("O4").Value = [double]"0.003431828163196655"
$ws.Range("P4").Value = [double]"128.8654118883123"
$ws.Range("Q4").Value = [double]"193.465830606327"

$ws.Range("A5").Value = "model_34_7_21"
$ws.Range("B5").Value = [double]"0.999988590433898"
$ws.Range("C5").Value = [double]"0.9990481780038731"
$ws.Range("D5").Value = [double]"0.9998364222709296"
$ws.Range("E5").Value = [double]"0.9999146950487101"
$ws.Range("F5").Value = [double]"0.999958254869512"
$ws.Range("G5").Value = [double]"1.065032913415572e-05"
$ws.Range("H5").Value = [double]"0.0008884840532275025"
$ws.Range("I5").Value = [double]"4.082536524128896e-05"
$ws.Range("J5").Value = [double]"2.452441349725371e-05"
$ws.Range("K5").Value = [double]"3.267488936927134e-05"
$ws.Range("L5").Value = [double]"0.0002421875183975704"
$ws.Range("M5").Value = [double]"0.003263484201609641"
$ws.Range("N5").Value = [double]"1.000009442399533"
$ws.Range("O5").Value = [double]"0.003402417480415774"
$ws.Range("P5").Value = [double]"128.8998395233416"
$ws.Range("Q5").Value = [double]"193.5002582413562"

$ws.Range("A6").Value = "model_34_7_20"
$ws.Range("B6").Value = [double]"0.9999888025139299"
$ws.Range("C6").Value = [double]"0.9990474102795875"
$ws.Range("D6").Value = [double]"0.9998402529405244"
$ws.Range("E6").Value = [double]"0.9999166652948901"
$ws.Range("F6").Value = [double]"0.9999592274214115"
$ws.Range("G6").Value = [double]"1.045236173361815e-05"
$ws.Range("H6").Value = [double]"0.000889200690148892"
$ws.Range("I6").Value = [double]"3.986931525687842e-05"
$ws.Range("J6").Value = [double]"2.395798527383127e-05"
$ws.Range("K6").Value = [double]"3.191365026535485e-05"
$ws.Range("L6").Value = [double]"0.0002430270512331329"
$ws.Range("M6").Value = [double]"0.003233011248606808"
$ws.Range("N6").Value = [double]"1.000009266885024"
$ws.Range("O6").Value = [double]"0.003370647230715901"
$ws.Range("P6").Value = [double]"128.9373652037426"
$ws.Range("Q6").Value = [double]"193.5377839217572"

$ws.Range("A7").Value = "model_34_7_19"
$ws.Range("B7").Value = [double]"0.9999890303050661"
$ws.Range("C7").Value = [double]"0.999046544695703"
$ws.Range("D7").Value = [double]"0.9998444664975048"
$ws.Range("E7").Value = [double]"0.9999188564246941"
$ws.Range("F7").Value = [double]"0.9999603018618133"
$ws.Range("G7").Value = [double]"1.023972870680745e-05"
$ws.Range("H7").Value = [double]"0.0008900086747104258"
$ws.Range("I7").Value = [double]"3.881770509168522e-05"
$ws.Range("J7").Value = [double]"2.332805497636143e-05"
$ws.Range("K7").Value = [double]"3.107266064935532e-05"
$ws.Range("L7").Value = [double]"0.0002439324772409027"
$ws.Range("M7").Value = [double]"0.003199957610157898"
$ws.Range("N7").Value = [double]"1.000009078368221"
$ws.Range("O7").Value = [double]"0.003336186430447757"
$ws.Range("P7").Value = [double]"128.9784708643594"
$ws.Range("Q7").Value = [double]"193.578889582374"

$ws.Range("A8").Value = "model_34_7_18"
$ws.Range("B8").Value = [double]"0.9999892713776307"
$ws.Range("C8").Value = [double]"0.9990455733925779"
$ws.Range("D8").Value = [double]"0.9998490784090912"
$ws.Range("E8").Value = [double]"0.9999212612546573"
$ws.Range("F8").Value = [double]"0.9999614784965621"
$ws.Range("G8").Value = [double]"1.001469804958e-05"
$ws.Range("H8").Value = [double]"0.0008909153435423234"
$ws.Range("I8").Value = [double]"3.766667447127529e-05"
$ws.Range("J8").Value = [double]"2.26366877870695e-05"
$ws.Range("K8").Value = [double]"3.01516811291724e-05"
$ws.Range("L8").Value = [double]"0.0002449313748567925"
$ws.Range("M8").Value = [double]"0.003164600772543039"
$ws.Range("N8").Value = [double]"1.000008878859892"
$ws.Range("O8").Value = [double]"0.003299324379056889"
$ws.Range("P8").Value = [double]"129.0229134782366"
$ws.Range("Q8").Value = [double]"193.6233321962512"

$ws.Range("A9").Value = "model_34_7_17"
$ws.Range("B9").Value = [double]"0.9999895244108192"
$ws.Range("C9").Value = [double]"0.9990444712399753"
$ws.Range("D9").Value = [double]"0.999854150171967"
$ws.Range("E9").Value = [double]"0.9999238668790216"
$ws.Range("F9").Value = [double]"0.9999627656007546"
$ws.Range("G9").Value = [double]"9.778502675025979e-06"
$ws.Range("H9").Value = [double]"0.000891944154618031"
$ws.Range("I9").Value = [double]"3.640087519043072e-05"
$ws.Range("J9").Value = [double]"2.188759399637138e-05"
$ws.Range("K9").Value = [double]"2.914423459340106e-05"
$ws.Range("L9").Value = [double]"0.00024602444330572"
$ws.Range("M9").Value = [double]"0.003127059749193478"
$ws.Range("N9").Value = [double]"1.000008669453115"
$ws.Range("O9").Value = [double]"0.003260185156622706"
$ws.Range("P9").Value = [double]"129.0706483727035"
$ws.Range("Q9").Value = [double]"193.6710670907181"

$ws.Range("A10").Value = "model_34_7_16"
$ws.Range("B10").Value = [double]"0.9999897884278663"
$ws.Range("C10").Value = [double]"0.9990432231181724"
$ws.Range("D10").Value = [double]"0.9998596937071356"
$ws.Range("E10").Value = [double]"0.9999267252487792"
$ws.Range("F10").Value = [double]"0.9999641743353551"
$ws.Range("G10").Value = [double]"9.532054350545965e-06"
$ws.Range("H10").Value = [double]"0.0008931092215348033"
$ws.Range("I10").Value = [double]"3.501733203165715e-05"
$ws.Range("J10").Value = [double]"2.106583815679865e-05"
$ws.Range("K10").Value = [double]"2.80415850942279e-05"
$ws.Range("L10").Value = [double]"0.0002472191529017118"
$ws.Range("M10").Value = [double]"0.003087402524865516"
$ws.Range("N10").Value = [double]"1.000008450956249"
$ws.Range("O10").Value = [double]"0.003218839642153332"
$ws.Range("P10").Value = [double]"129.1217005936807"
$ws.Range("Q10").Value = [double]"193.7221193116953"

$ws.Range("A11").Value = "model_34_7_15"
$ws.Range("B11").Value = [double]"0.9999900609555421"
$ws.Range("C11").Value = [double]"0.9990418126100765"
$ws.Range("D11").Value = [double]"0.9998657437228848"
$ws.Range("E11").Value = [double]"0.9999298623799769"
$ws.Range("F11").Value = [double]"0.9999657150135894"
$ws.Range("G11").Value = [double]"9.277661727771472e-06"
$ws.Range("H11").Value = [double]"0.0008944258689282069"
$ws.Range("I11").Value = [double]"3.350738257773665e-05"
$ws.Range("J11").Value = [double]"2.016394088679385e-05"
$ws.Range("K11").Value = [double]"2.683566022899289e-05"
$ws.Range("L11").Value = [double]"0.0002485085863572947"
$ws.Range("M11").Value = [double]"0.003045925430435136"
$ws.Range("N11").Value = [double]"1.000008225416103"
$ws.Range("O11").Value = [double]"0.003175596781943629"
$ws.Range("P11").Value = [double]"129.1758020238117"
$ws.Range("Q11").Value = [double]"193.7762207418263"

$ws.Range("A12").Value = "model_34_7_14"
$ws.Range("B12").Value = [double]"0.999990335449617"
$ws.Range("C12").Value = [double]"0.9990402117571732"
$ws.Range("D12").Value = [double]"0.9998723181397355"
$ws.Range("E12").Value = [double]"0.999933256114513"
$ws.Range("F12").Value = [double]"0.9999673864174582"
$ws.Range("G12").Value = [double]"9.021433557751277e-06"
$ws.Range("H12").Value = [double]"0.0008959201948441405"
$ws.Range("I12").Value = [double]"3.186655426510048e-05"
$ws.Range("J12").Value = [double]"1.918827244312743e-05"
$ws.Range("K12").Value = [double]"2.552741335411396e-05"
$ws.Range("L12").Value = [double]"0.0002500413838968197"
$ws.Range("M12").Value = [double]"0.003003570135314186"
$ws.Range("N12").Value = [double]"1.000007998248593"
$ws.Range("O12").Value = [double]"0.003131438334221831"
$ws.Range("P12").Value = [double]"129.2318146110068"
$ws.Range("Q12").Value = [double]"193.8322333290214"

$ws.Range("A13").Value = "model_34_7_13"
$ws.Range("B13").Value = [double]"0.9999906074474645"
$ws.Range("C13").Value = [double]"0.9990383858632854"
$ws.Range("D13").Value = [double]"0.999879456357406"
$ws.Range("E13").Value = [double]"0.9999369349104534"
$ws.Range("F13").Value = [double]"0.9999692013372277"
$ws.Range("G13").Value = [double]"8.767535506502441e-06"
$ws.Range("H13").Value = [double]"0.0008976245866408275"
$ws.Range("I13").Value = [double]"3.008501379974381e-05"
$ws.Range("J13").Value = [double]"1.8130651385357e-05"
$ws.Range("K13").Value = [double]"2.410683322922688e-05"
$ws.Range("L13").Value = [double]"0.0002518402799276048"
$ws.Range("M13").Value = [double]"0.00296100244959413"
$ws.Range("N13").Value = [double]"1.000007773146926"
$ws.Range("O13").Value = [double]"0.003087058453993416"
$ws.Range("P13").Value = [double]"129.2889096103077"
$ws.Range("Q13").Value = [double]"193.8893283283223"

$ws.Range("A14").Value = "model_34_7_12"
$ws.Range("B14").Value = [double]"0.9999908679996334"
$ws.Range("C14").Value = [double]"0.9990363069821123"
$ws.Range("D14").Value = [double]"0.9998871810645853"
$ws.Range("E14").Value = [double]"0.9999409204232388"
$ws.Range("F14").Value = [double]"0.9999711635348605"
$ws.Range("G14").Value = [double]"8.524321493704776e-06"
$ws.Range("H14").Value = [double]"0.0008995651309636177"
$ws.Range("I14").Value = [double]"2.815709858924234e-05"
$ws.Range("J14").Value = [double]"1.698485196728289e-05"
$ws.Range("K14").Value = [double]"2.257097527826261e-05"
$ws.Range("L14").Value = [double]"0.0002537662729709331"
$ws.Range("M14").Value = [double]"0.00291964406969493"
$ws.Range("N14").Value = [double]"1.000007557517545"
$ws.Range("O14").Value = [double]"0.003043939362238258"
$ws.Range("P14").Value = [double]"129.3451742562868"
$ws.Range("Q14").Value = [double]"193.9455929743015"

$ws.Range("A15").Value = "model_34_7_11"
$ws.Range("B15").Value = [double]"0.9999911059230286"
$ws.Range("C15").Value = [double]"0.9990339318336829"
$ws.Range("D15").Value = [double]"0.9998955084203665"
$ws.Range("E15").Value = [double]"0.9999451815804352"
$ws.Range("F15").Value = [double]"0.9999732736971989"
$ws.Range("G15").Value = [double]"8.302230447887982e-06"
$ws.Range("H15").Value = [double]"0.0009017822277654683"
$ws.Range("I15").Value = [double]"2.607877568311247e-05"
$ws.Range("J15").Value = [double]"1.575980723679899e-05"
$ws.Range("K15").Value = [double]"2.091930189380133e-05"
$ws.Range("L15").Value = [double]"0.0002558281342741295"
$ws.Range("M15").Value = [double]"0.00288135913205695"
$ws.Range("N15").Value = [double]"1.000007360615425"
$ws.Range("O15").Value = [double]"0.003004024555544283"
$ws.Range("P15").Value = [double]"129.3979727012032"
$ws.Range("Q15").Value = [double]"193.9983914192178"

$ws.Range("A16").Value = "model_34_7_10"
$ws.Range("B16").Value = [double]"0.9999913056622838"
$ws.Range("C16").Value = [double]"0.9990312109855078"
$ws.Range("D16").Value = [double]"0.9999043946547529"
$ws.Range("E16").Value = [double]"0.9999497469048757"
$ws.Range("F16").Value = [double]"0.9999755305781715"
$ws.Range("G16").Value = [double]"8.115782620811859e-06"
$ws.Range("H16").Value = [double]"0.0009043220201054882"
$ws.Range("I16").Value = [double]"2.386096909959929e-05"
$ws.Range("J16").Value = [double]"1.4447317133528e-05"
$ws.Range("K16").Value = [double]"1.915278840503909e-05"
$ws.Range("L16").Value = [double]"0.00025802901860762"
$ws.Range("M16").Value = [double]"0.002848821268667422"
$ws.Range("N16").Value = [double]"1.000007195313972"
$ws.Range("O16").Value = [double]"0.002970101487947599"
$ws.Range("P16").Value = [double]"129.443399840779"
$ws.Range("Q16").Value = [double]"194.0438185587936"

$ws.Range("A17").Value = "model_34_7_9"
$ws.Range("B17").Value = [double]"0.9999914443557965"
$ws.Range("C17").Value = [double]"0.9990280856065558"
$ws.Range("D17").Value = [double]"0.9999138644450088"
$ws.Range("E17").Value = [double]"0.9999545454054328"
$ws.Range("F17").Value = [double]"0.9999779198423864"
$ws.Range("G17").Value = [double]"7.986318314664022e-06"
$ws.Range("H17").Value = [double]"0.0009072394241689308"
$ws.Range("I17").Value = [double]"2.14975199421055e-05"
$ws.Range("J17").Value = [double]"1.306779097415579e-05"
$ws.Range("K17").Value = [double]"1.728265545813064e-05"
$ws.Range("L17").Value = [double]"0.0002603829937706274"
$ws.Range("M17").Value = [double]"0.002826007486661"
$ws.Range("N17").Value = [double]"1.000007080533134"
$ws.Range("O17").Value = [double]"0.002946316476010125"
$ws.Range("P17").Value = [double]"129.4755613820624"
$ws.Range("Q17").Value = [double]"194.075980100077"

$ws.Range("A18").Value = "model_34_7_8"
$ws.Range("B18").Value = [double]"0.9999914968214574"
$ws.Range("C18").Value = [double]"0.9990244975983417"
$ws.Range("D18").Value = [double]"0.9999238303740725"
$ws.Range("E18").Value = [double]"0.9999595577998942"
$ws.Range("F18").Value = [double]"0.9999804292175114"
$ws.Range("G18").Value = [double]"7.937343923098278e-06"
$ws.Range("H18").Value = [double]"0.000910588672341474"
$ws.Range("I18").Value = [double]"1.901024556614389e-05"
$ws.Range("J18").Value = [double]"1.162677222290396e-05"
$ws.Range("K18").Value = [double]"1.531850889452392e-05"
$ws.Range("L18").Value = [double]"0.0002628666422337683"
$ws.Range("M18").Value = [double]"0.002817329218089054"
$ws.Range("N18").Value = [double]"1.000007037113277"
$ws.Range("O18").Value = [double]"0.002937268755578578"
$ws.Range("P18").Value = [double]"129.487863714341"
$ws.Range("Q18").Value = [double]"194.0882824323556"

$ws.Range("A19").Value = "model_34_7_7"
$ws.Range("B19").Value = [double]"0.999991430926731"
$ws.Range("C19").Value = [double]"0.9990203570852506"
$ws.Range("D19").Value = [double]"0.9999342444036029"
$ws.Range("E19").Value = [double]"0.9999647448654566"
$ws.Range("F19").Value = [double]"0.9999830420847715"
$ws.Range("G19").Value = [double]"7.998853757775142e-06"
$ws.Range("H19").Value = [double]"0.0009144536595644727"
$ws.Range("I19").Value = [double]"1.64111352738944e-05"
$ws.Range("J19").Value = [double]"1.013553708630995e-05"
$ws.Range("K19").Value = [double]"1.327335661781834e-05"
$ws.Range("L19").Value = [double]"0.0002671980252461115"
$ws.Range("M19").Value = [double]"0.002828224488574969"
$ws.Range("N19").Value = [double]"1.000007091646843"
$ws.Range("O19").Value = [double]"0.002948627860285397"
$ws.Range("P19").Value = [double]"129.4724246136563"
$ws.Range("Q19").Value = [double]"194.0728433316709"

$ws.Range("A20").Value = "model_34_7_6"
$ws.Range("B20").Value = [double]"0.9999912008203932"
$ws.Range("C20").Value = [double]"0.9990155666544931"
$ws.Range("D20").Value = [double]"0.9999449689092837"
$ws.Range("E20").Value = [double]"0.9999699869315161"
$ws.Range("F20").Value = [double]"0.9999857146039168"
$ws.Range("G20").Value = [double]"8.213647923634781e-06"
$ws.Range("H20").Value = [double]"0.0009189253164010161"
$ws.Range("I20").Value = [double]"1.373453703559484e-05"
$ws.Range("J20").Value = [double]"8.628489796820938e-06"
$ws.Range("K20").Value = [double]"1.118151341620789e-05"
$ws.Range("L20").Value = [double]"0.0002730940878298339"
$ws.Range("M20").Value = [double]"0.002865946252747037"
$ws.Range("N20").Value = [double]"1.000007282079675"
$ws.Range("O20").Value = [double]"0.002987955518053087"
$ws.Range("P20").Value = [double]"129.41942681262"
$ws.Range("Q20").Value = [double]"194.0198455306346"

$ws.Range("A21").Value = "model_34_7_5"
$ws.Range("B21").Value = [double]"0.9999907453202869"
$ws.Range("C21").Value = [double]"0.9990100153968007"
$ws.Range("D21").Value = [double]"0.9999557969047691"
$ws.Range("E21").Value = [double]"0.9999751402168872"
$ws.Range("F21").Value = [double]"0.9999883872918462"
$ws.Range("G21").Value = [double]"8.638837278693549e-06"
$ws.Range("H21").Value = [double]"0.0009241071717847739"
$ws.Range("I21").Value = [double]"1.103211004243721e-05"
$ws.Range("J21").Value = [double]"7.146966164253599e-06"
$ws.Range("K21").Value = [double]"9.089538103345406e-06"
$ws.Range("L21").Value = [double]"0.0002794696057635864"
$ws.Range("M21").Value = [double]"0.00293918990177456"
$ws.Range("N21").Value = [double]"1.00000765904528"
$ws.Range("O21").Value = [double]"0.003064317300854968"
$ws.Range("P21").Value = [double]"129.318485116859"
$ws.Range("Q21").Value = [double]"193.9189038348736"

$ws.Range("A22").Value = "model_34_7_4"
$ws.Range("B22").Value = [double]"0.9999899957781573"
$ws.Range("C22").Value = [double]"0.9990035531607461"
$ws.Range("D22").Value = [double]"0.9999664954289857"
$ws.Range("E22").Value = [double]"0.9999799862864572"
$ws.Range("F22").Value = [double]"0.9999909829187343"
$ws.Range("G22").Value = [double]"9.338501955579515e-06"
$ws.Range("H22").Value = [double]"0.0009301393854823805"
$ws.Range("I22").Value = [double]"8.361996200111301e-06"
$ws.Range("J22").Value = [double]"5.753764337450641e-06"
$ws.Range("K22").Value = [double]"7.057880268780972e-06"
$ws.Range("L22").Value = [double]"0.0002862644114600338"
$ws.Range("M22").Value = [double]"0.003055896260604982"
$ws.Range("N22").Value = [double]"1.000008279356008"
$ws.Range("O22").Value = [double]"0.003185992090996269"
$ws.Range("P22").Value = [double]"129.1627294175646"
$ws.Range("Q22").Value = [double]"193.7631481355793"

$ws.Range("A23").Value = "model_34_7_3"
$ws.Range("B23").Value = [double]"0.9999888558952"
$ws.Range("C23").Value = [double]"0.9989960282896652"
$ws.Range("D23").Value = [double]"0.9999766561243671"
$ws.Range("E23").Value = [double]"0.9999842830330034"
$ws.Range("F23").Value = [double]"0.9999933919180041"
$ws.Range("G23").Value = [double]"1.040253266110057e-05"
$ws.Range("H23").Value = [double]"0.0009371635223327338"
$ws.Range("I23").Value = [double]"5.826112480421737e-06"
$ws.Range("J23").Value = [double]"4.518487985994099e-06"
$ws.Range("K23").Value = [double]"5.172300233207919e-06"
$ws.Range("L23").Value = [double]"0.0002935451838447891"
$ws.Range("M23").Value = [double]"0.003225295747850199"
$ws.Range("N23").Value = [double]"1.000009222707421"
$ws.Range("O23").Value = [double]"0.003362603265118795"
$ws.Range("P23").Value = [double]"128.9469225127171"
$ws.Range("Q23").Value = [double]"193.5473412307317"

$ws.Range("A24").Value = "model_34_7_2"
$ws.Range("B24").Value = [double]"0.9999872047676039"
$ws.Range("C24").Value = [double]"0.9989872142343635"
$ws.Range("D24").Value = [double]"0.9999858136238371"
$ws.Range("E24").Value = [double]"0.9999876190858719"
$ws.Range("F24").Value = [double]"0.9999954645492211"
$ws.Range("G24").Value = [double]"1.194378779592619e-05"
$ws.Range("H24").Value = [double]"0.0009453910560644908"
$ws.Range("I24").Value = [double]"3.540604161658108e-06"
$ws.Range("J24").Value = [double]"3.559402507853633e-06"
$ws.Range("K24").Value = [double]"3.550003334755871e-06"
$ws.Range("L24").Value = [double]"0.0003014294403009376"
$ws.Range("M24").Value = [double]"0.003455978558371881"
$ws.Range("N24").Value = [double]"1.000010589157845"
$ws.Range("O24").Value = [double]"0.003603106720463634"
$ws.Range("P24").Value = [double]"128.6705985289375"
$ws.Range("Q24").Value = [double]"193.2710172469521"

$ws.Range("A25").Value = "model_34_7_1"
$ws.Range("B25").Value = [double]"0.9999848993784193"
$ws.Range("C25").Value = [double]"0.9989769167964646"
$ws.Range("D25").Value = [double]"0.9999932834956616"
$ws.Range("E25").Value = [double]"0.9999895177355445"
$ws.Range("F25").Value = [double]"0.9999970041473197"
$ws.Range("G25").Value = [double]"1.409576740494591e-05"
$ws.Range("H25").Value = [double]"0.0009550032623379442"
$ws.Range("I25").Value = [double]"1.676290191327166e-06"
$ws.Range("J25").Value = [double]"3.013557642405273e-06"
$ws.Range("K25").Value = [double]"2.344923916866219e-06"
$ws.Range("L25").Value = [double]"0.0003099040243903889"
$ws.Range("M25").Value = [double]"0.003754433033754353"
$ws.Range("N25").Value = [double]"1.000012497066136"
$ws.Range("O25").Value = [double]"0.003914267020749073"
$ws.Range("P25").Value = [double]"128.3392719793816"
$ws.Range("Q25").Value = [double]"192.9396906973962"

$ws.Range("A26").Value = "model_34_7_0"
$ws.Range("B26").Value = [double]"0.9999817543205211"
$ws.Range("C26").Value = [double]"0.9989648146925866"
$ws.Range("D26").Value = [double]"0.9999982858833012"
$ws.Range("E26").Value = [double]"0.9999893451083597"
$ws.Range("F26").Value = [double]"0.9999977699708703"
$ws.Range("G26").Value = [double]"1.703154090081487e-05"
$ws.Range("H26").Value = [double]"0.000966300045087021"
$ws.Range("I26").Value = [double]"4.278054273792097e-07"
$ws.Range("J26").Value = [double]"3.063186420078487e-06"
$ws.Range("K26").Value = [double]"1.745495923728848e-06"
$ws.Range("L26").Value = [double]"0.0003190469745343864"
$ws.Range("M26").Value = [double]"0.004126928749180784"
$ws.Range("N26").Value = [double]"1.000015099872672"
$ws.Range("O26").Value = [double]"0.004302620649953638"
$ws.Range("P26").Value = [double]"127.9608871716636"
$ws.Range("Q26").Value = [double]"192.5613058896782"
